$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of the existing "Unique Paths" row (row 6) down onto the
# new row 7 so the new entry matches the look of the rest of the table
# (gray-highlighted cells, left/top aligned + wrapped text, etc.)
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new "Unique Paths II" row
$ws.Range("A7").Value = 63
$ws.Range("B7").Value = "LC/CN/GFG"
$ws.Range("C7").Value = "Unique Paths II"
$ws.Range("D7").Value = "Java"
$ws.Range("E7").Value = "DP(Tabulation+space optimization)"
$ws.Range("F7").Value = "https://leetcode.com/problems/unique-paths/description/"

# Update the active view: scroll back so column A is visible again and select C7
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select()
